$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: D: '26.743.30' -> '26.754.58', E: '  +0.28%  ' -> '  +0.30%  '
$ws.Range("D2").Value = "'26.754.58"
$ws.Range("E2").Value = "  +0.30%  "

# Row 3: D: '1.603.32' -> '1.604.04', E: '  +0.38%  ' -> '  +0.40%  '
$ws.Range("D3").Value = "'1.604.04"
$ws.Range("E3").Value = "  +0.40%  "

# Row 4: E: '  +0.22%  ' -> '  +0.20%  '
$ws.Range("E4").Value = "  +0.20%  "

# Row 5: D: '211.82' -> '211.87', E: '  +0.08%  ' -> '  +0.16%  '
$ws.Range("D5").Value = "'211.87"
$ws.Range("E5").Value = "  +0.16%  "

# Row 6: E: '  +0.01%  ' -> '  +0.13%  '
$ws.Range("E6").Value = "  +0.13%  "

# Row 7: E: '  +0.23%  ' -> '  +0.21%  '
$ws.Range("E7").Value = "  +0.21%  "

# Row 8: E: '  -0.16%  ' -> '  -0.06%  '
$ws.Range("E8").Value = "  -0.06%  "

# Row 9: E: '  -0.09%  ' -> '  +0.04%  '
$ws.Range("E9").Value = "  +0.04%  "

# Row 10: E: '  +0.56%  ' -> '  +0.64%  '
$ws.Range("E10").Value = "  +0.64%  "

# Row 11: E: '  +0.70%  ' -> '  +0.80%  '
$ws.Range("E11").Value = "  +0.80%  "

# Row 12: E: '  +0.41%  ' -> '  +0.39%  '
$ws.Range("E12").Value = "  +0.39%  "

# Row 13: D: '1.603.91' -> '1.607.40', E: '  +0.32%  ' -> '  -0.35%  '
$ws.Range("D13").Value = "'1.607.40"
$ws.Range("E13").Value = "  -0.35%  "

# Row 14: E: '  +0.88%  ' -> '  +0.97%  '
$ws.Range("E14").Value = "  +0.97%  "

# Row 15: E: '  +0.33%  ' -> '  +0.43%  '
$ws.Range("E15").Value = "  +0.43%  "

# Row 16: D: '65.07' -> '65.10', E: '  -0.03%  ' -> '  +0.04%  '
$ws.Range("D16").Value = "'65.10"
$ws.Range("E16").Value = "  +0.04%  "

# Row 17: D: '0.0₃0741' -> '0.0₃0742', E: '  -1.56%  ' -> '  -1.23%  '
$ws.Range("D17").Value = "'0.0₃0742"
$ws.Range("E17").Value = "  -1.23%  "

# Row 18: D: '7.20' -> '7.21', E: '  +1.82%  ' -> '  +1.97%  '
$ws.Range("D18").Value = "'7.21"
$ws.Range("E18").Value = "  +1.97%  "

# Row 19: B: 'Dai' -> 'BitcoinCash', C: 'https://coinranking.com/coin/MoTuySvg7+dai-dai' -> 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', D: '1.01' -> '209.75', E: '  +0.22%  ' -> '  -0.20%  '
$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").Value = "'209.75"
$ws.Range("E19").Value = "  -0.20%  "

# Row 20: B: 'BitcoinCash' -> 'Dai', C: 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch' -> 'https://coinranking.com/coin/MoTuySvg7+dai-dai', D: '209.10' -> '1.01', E: '  -0.51%  ' -> '  +0.23%  '
$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D20").Value = "'1.01"
$ws.Range("E20").Value = "  +0.23%  "

# Row 21: D: '4.30' -> '4.31', E: '  +0.23%  ' -> '  +0.31%  '
$ws.Range("D21").Value = "'4.31"
$ws.Range("E21").Value = "  +0.31%  "

# Row 22: E: '  -4.64%  ' -> '  -4.65%  '
$ws.Range("E22").Value = "  -4.65%  "

# Row 23: D: '9.04' -> '9.05', E: '  +0.77%  ' -> '  +0.89%  '
$ws.Range("D23").Value = "'9.05"
$ws.Range("E23").Value = "  +0.89%  "

# Row 24: D: '143.79' -> '143.81', E: '  +0.35%  ' -> '  +0.37%  '
$ws.Range("D24").Value = "'143.81"
$ws.Range("E24").Value = "  +0.37%  "

# Row 25: E: '  +0.20%  ' -> '  +0.22%  '
$ws.Range("E25").Value = "  +0.22%  "

# Row 26: D: '7.10' -> '7.13', E: '  -0.15%  ' -> '  +0.12%  '
$ws.Range("D26").Value = "'7.13"
$ws.Range("E26").Value = "  +0.12%  "

# Row 27: E: '  -0.20%  ' -> '  -0.15%  '
$ws.Range("E27").Value = "  -0.15%  "

# Row 28: D: '15.36' -> '15.37', E: '  +0.10%  ' -> '  +0.09%  '
$ws.Range("D28").Value = "'15.37"
$ws.Range("E28").Value = "  +0.09%  "

# Row 30: E: '  -0.03%  ' -> '  +0.14%  '
$ws.Range("E30").Value = "  +0.14%  "

# Row 31: E: '  +0.90%  ' -> '  +0.84%  '
$ws.Range("E31").Value = "  +0.84%  "

# Row 32: E: '  +0.46%  ' -> '  +0.35%  '
$ws.Range("E32").Value = "  +0.35%  "

# Row 33: D: '1.287.79' -> '1.287.56', E: '  -0.02%  ' -> '  +0.00%  '
$ws.Range("D33").Value = "'1.287.56"
$ws.Range("E33").Value = "  +0.00%  "

# Row 34: E: '  +1.38%  ' -> '  +1.35%  '
$ws.Range("E34").Value = "  +1.35%  "

# Row 35: E: '  +16.90%  ' -> '  +16.18%  '
$ws.Range("E35").Value = "  +16.18%  "

# Row 36: E: '  +0.17%  ' -> '  +0.19%  '
$ws.Range("E36").Value = "  +0.19%  "

# Row 37: D: '0.588' -> '0.589', E: '  -5.16%  ' -> '  -4.99%  '
$ws.Range("D37").Value = "'0.589"
$ws.Range("E37").Value = "  -4.99%  "

# Row 38: E: '  -0.87%  ' -> '  -0.85%  '
$ws.Range("E38").Value = "  -0.85%  "

# Row 39: D: '0.826' -> '0.827', E: '  -0.20%  ' -> '  +0.05%  '
$ws.Range("D39").Value = "'0.827"
$ws.Range("E39").Value = "  +0.05%  "

# Row 40: E: '  -0.22%  ' -> '  -0.12%  '
$ws.Range("E40").Value = "  -0.12%  "

# Row 41: E: '  +0.11%  ' -> '  +0.32%  '
$ws.Range("E41").Value = "  +0.32%  "

# Row 42: E: '  -0.44%  ' -> '  -0.49%  '
$ws.Range("E42").Value = "  -0.49%  "

# Row 43: D: '62.65' -> '62.71', E: '  -1.09%  ' -> '  -0.92%  '
$ws.Range("D43").Value = "'62.71"
$ws.Range("E43").Value = "  -0.92%  "

# Row 44: D: '1.740.55' -> '1.740.31', E: '  +0.39%  ' -> '  +0.35%  '
$ws.Range("D44").Value = "'1.740.31"
$ws.Range("E44").Value = "  +0.35%  "

# Row 45: D: '90.42' -> '90.39', E: '  -0.94%  ' -> '  -1.04%  '
$ws.Range("D45").Value = "'90.39"
$ws.Range("E45").Value = "  -1.04%  "

# Row 46: E: '  +0.04%  ' -> '  +0.15%  '
$ws.Range("E46").Value = "  +0.15%  "

# Row 47: B: 'BabyDogeCoin' -> 'Algorand', C: 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge' -> 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo', D: '0.0₆0103' -> '0.102', E: '  -2.67%  ' -> '  +0.80%  '
$ws.Range("B47").Value = "Algorand"
$ws.Range("C47").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D47").Value = "'0.102"
$ws.Range("E47").Value = "  +0.80%  "

# Row 48: B: 'Algorand' -> 'Cronos', C: 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo' -> 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro', D: '0.102' -> '0.0513', E: '  +0.66%  ' -> '  +0.75%  '
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").Value = "'0.0513"
$ws.Range("E48").Value = "  +0.75%  "

# Row 49: B: 'Cronos' -> 'EnergySwap', C: 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro' -> 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens', D: '0.0513' -> '7.55', E: '  +0.69%  ' -> '  +3.19%  '
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'7.55"
$ws.Range("E49").Value = "  +3.19%  "

# Row 50: B: 'EnergySwap' -> 'USDD', C: 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens' -> 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd', D: '7.54' -> '1.00', E: '  +3.00%  ' -> '  +0.15%  '
$ws.Range("B50").Value = "USDD"
$ws.Range("C50").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D50").Value = "'1.00"
$ws.Range("E50").Value = "  +0.15%  "

# Row 51: B: 'USDD' -> 'Mantle', C: 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd' -> 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt', D: '1.00' -> '0.399', E: '  +0.14%  ' -> '  +1.88%  '
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").Value = "'0.399"
$ws.Range("E51").Value = "  +1.88%  "
